$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.698.52"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.450.54"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.00"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.74"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.486"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.03"
$ws.Range("E9").Value = "  +4.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.412"
$ws.Range("E11").Value = "  +2.77%  "
$ws.Range("D12").Value = "4.034.33"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.36"
$ws.Range("E14").Value = "  -5.36%  "
$ws.Range("D15").Value = "3.451.42"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "62.733.97"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.42"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.67"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.01"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.03"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.566"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.28"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "3.583.24"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.183"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.01"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.25"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.21"
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.95"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.06"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").Value = "3.483.94"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0775"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.69"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.39"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").Value = "2.573.79"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.25"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.58"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("E51").Value = "  +0.07%  "
